$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.895.71"
$ws.Range("E2").Value = "  +2.49%  "

# Row 3
$ws.Range("D3").Value = "3.568.42"
$ws.Range("E3").Value = "  +1.48%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'583.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.17%  "

# Row 6
$ws.Range("D6").Value = "'186.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.92%  "

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.84%  "

# Row 8
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.557.88"
$ws.Range("E8").Value = "  +1.31%  "

# Row 9
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("E10").Value = "  +17.16%  "

# Row 11
$ws.Range("E11").Value = "  +2.67%  "

# Row 12
$ws.Range("D12").Value = "'54.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.39%  "

# Row 13
$ws.Range("D13").Value = "'0.0000319"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.91%  "

# Row 14
$ws.Range("D14").Value = "'9.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.47%  "

# Row 15
$ws.Range("D15").Value = "4.133.92"

# Row 16
$ws.Range("D16").Value = "70.906.44"
$ws.Range("E16").Value = "  +2.60%  "

# Row 17
$ws.Range("D17").Value = "'19.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "

# Row 18
$ws.Range("D18").Value = "3.592.63"
$ws.Range("E18").Value = "  +2.03%  "

# Row 19
$ws.Range("D19").Value = "'580.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.61%  "

# Row 20
$ws.Range("D20").Value = "'12.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.36%  "

# Row 21
$ws.Range("E21").Value = "  +0.78%  "

# Row 22
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.18%  "

# Row 23
$ws.Range("E23").Value = "  -14.34%  "

# Row 24
$ws.Range("D24").Value = "'5.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "

# Row 25
$ws.Range("D25").Value = "'4.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.21%  "

# Row 26
$ws.Range("D26").Value = "'95.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.47%  "

# Row 27
$ws.Range("E27").Value = "  +2.63%  "

# Row 28
$ws.Range("D28").Value = "'2.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.24%  "

# Row 29
$ws.Range("D29").Value = "'9.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.36%  "

# Row 30
$ws.Range("D30").Value = "'32.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.95%  "

# Row 31
$ws.Range("D31").Value = "'7.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "

# Row 32
$ws.Range("E32").Value = "  -3.30%  "

# Row 33
$ws.Range("D33").Value = "'0.116"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.26%  "

# Row 34
$ws.Range("D34").Value = "'63.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.55%  "

# Row 35
$ws.Range("D35").Value = "'3.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.02%  "

# Row 36
$ws.Range("D36").Value = "'549.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.40%  "

# Row 37
$ws.Range("D37").Value = "'0.416"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.57%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0802"
$ws.Range("E38").Value = "  +5.23%  "

# Row 39
$ws.Range("D39").Value = "'37.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.21%  "

# Row 40
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("D41").Value = "'3.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.89%  "

# Row 42
$ws.Range("D42").Value = "3.564.19"
$ws.Range("E42").Value = "  +11.85%  "

# Row 43
$ws.Range("E43").Value = "  +2.78%  "

# Row 44
$ws.Range("D44").Value = "'3.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.85%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.21%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0447"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.34%  "

# Row 47
$ws.Range("D47").Value = "'2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "

# Row 48
$ws.Range("D48").Value = "'9.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.91%  "

# Row 49
$ws.Range("D49").Value = "'0.138"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.00%  "

# Row 50
$ws.Range("D50").Value = "'1.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.47%  "

# Row 51
$ws.Range("E51").Value = "  +0.06%  "
